$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.48"
$ws.Range("E2").Value = "'1.35%"
$ws.Range("G2").Value = "'2"
$ws.Range("D3").Value = "'29.54"
$ws.Range("E3").Value = "'0.73%"
$ws.Range("G3").Value = "'2"
$ws.Range("D4").Value = "'5.154"
$ws.Range("E4").Value = "'1.19%"
$ws.Range("G4").Value = "'2"
$ws.Range("D5").Value = "'0.05747"
$ws.Range("E5").Value = "'1.80%"
$ws.Range("G5").Value = "'2"
$ws.Range("D6").Value = "'6.563"
$ws.Range("E6").Value = "'0.94%"
$ws.Range("G6").Value = "'2"
$ws.Range("D7").Value = "'0.8592"
$ws.Range("E7").Value = "'4.75%"
$ws.Range("G7").Value = "'2"
$ws.Range("D8").Value = "'0.8542"
$ws.Range("E8").Value = "'-0.02%"
$ws.Range("G8").Value = "'2"
$ws.Range("D9").Value = "'0.1361"
$ws.Range("E9").Value = "'2.40%"
$ws.Range("G9").Value = "'2"
$ws.Range("D10").Value = "'0.07020"
$ws.Range("E10").Value = "'1.34%"
$ws.Range("G10").Value = "'2"
$ws.Range("D11").Value = "'0.03035"
$ws.Range("E11").Value = "'5.81%"
$ws.Range("G11").Value = "'2"
$ws.Range("D12").Value = "'0.09361"
$ws.Range("E12").Value = "'-0.07%"
$ws.Range("G12").Value = "'2"
$ws.Range("D13").Value = "'0.001530"
$ws.Range("E13").Value = "'1.40%"
$ws.Range("G13").Value = "'2"
$ws.Range("D14").Value = "'0.0005980"
$ws.Range("E14").Value = "'-94.02%"
$ws.Range("G14").Value = "'2"
$ws.Range("D15").Value = "'0.006094"
$ws.Range("E15").Value = "'0.99%"
$ws.Range("G15").Value = "'2"
$ws.Range("D16").Value = "'3.502"
$ws.Range("E16").Value = "'-0.69%"
$ws.Range("G16").Value = "'2"
$ws.Range("D17").Value = "'3.106"
$ws.Range("E17").Value = "'2.97%"
$ws.Range("G17").Value = "'2"
$ws.Range("D18").Value = "'2.217"
$ws.Range("E18").Value = "'-0.06%"
$ws.Range("G18").Value = "'2"
$ws.Range("D19").Value = "'0.3198"
$ws.Range("E19").Value = "'1.50%"
$ws.Range("G19").Value = "'2"
$ws.Range("D20").Value = "'0.03271"
$ws.Range("E20").Value = "'1.84%"
$ws.Range("G20").Value = "'2"
$ws.Range("D21").Value = "'0.1282"
$ws.Range("E21").Value = "'0.72%"
$ws.Range("G21").Value = "'2"
$ws.Range("D22").Value = "'3.558"
$ws.Range("E22").Value = "'-1.53%"
$ws.Range("G22").Value = "'2"
$ws.Range("D23").Value = "'0.04143"
$ws.Range("E23").Value = "'0.10%"
$ws.Range("G23").Value = "'2"
$ws.Range("E24").Value = "'0.39%"
$ws.Range("G24").Value = "'2"
$ws.Range("E25").Value = "'1.20%"
$ws.Range("G25").Value = "'2"
$ws.Range("D26").Value = "'0.004129"
$ws.Range("E26").Value = "'-7.27%"
$ws.Range("G26").Value = "'2"
$ws.Range("G27").Value = "'2"
$ws.Range("E28").Value = "'3.00%"
$ws.Range("G28").Value = "'2"
$ws.Range("G29").Value = "'2"
$ws.Range("G30").Value = "'2"
$ws.Range("G31").Value = "'2"
$ws.Range("G32").Value = "'2"
$ws.Range("G33").Value = "'2"
$ws.Range("G34").Value = "'2"
$ws.Range("G35").Value = "'2"
$ws.Range("G36").Value = "'2"
$ws.Range("G37").Value = "'2"
$ws.Range("G38").Value = "'2"
$ws.Range("G39").Value = "'2"
$ws.Range("D40").Value = "'0.03726"
$ws.Range("E40").Value = "'0.69%"
$ws.Range("G40").Value = "'2"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.005869"
$ws.Range("E41").Value = "'-0.21%"
$ws.Range("G41").Value = "'2"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1068"
$ws.Range("E42").Value = "'1.42%"
$ws.Range("G42").Value = "'2"
$ws.Range("D43").Value = "'0.002200"
$ws.Range("E43").Value = "'8.56%"
$ws.Range("G43").Value = "'2"
$ws.Range("D44").Value = "'0.009314"
$ws.Range("E44").Value = "'-1.38%"
$ws.Range("G44").Value = "'2"
$ws.Range("D45").Value = "'0.00005273"
$ws.Range("E45").Value = "'3.42%"
$ws.Range("G45").Value = "'2"
$ws.Range("G46").Value = "'2"
$ws.Range("D47").Value = "'0.05800"
$ws.Range("E47").Value = "'-51.68%"
$ws.Range("G47").Value = "'2"
$ws.Range("D48").Value = "'0.002461"
$ws.Range("E48").Value = "'-2.37%"
$ws.Range("G48").Value = "'2"
$ws.Range("G49").Value = "'2"
$ws.Range("G50").Value = "'2"
$ws.Range("G51").Value = "'2"
